# Add a new "Serviced by " column (O) to the "Card5" sheet, and tidy up the
# "Correction " header (drop its trailing space) and backfill the existing
# "nan" placeholder pattern for the newly widened data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card5")

# 1) Header row: fix "Correction " -> "Correction" and add "Serviced by "
$ws.Range("N1").Value = "Correction"

# Clone the header's style (bold, centered, bordered) onto the new column's
# header cell before giving it its text.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Serviced by "

# 2) Data rows 2-13: the previously-blank "Correction" cells become the
# sheet-wide "nan" placeholder, matching every other data column.
$ws.Range("N2:N13").Value = "nan"

# 3) Data rows 2-13: materialize the new, still-blank "Serviced by " cells
# so the sheet's used range/dimension extends through column O.
$ws.Range("O2:O13").Font.Bold = $false

Write-Output "Card5 updated: N1/O1 headers, N2:N13 -> nan, O2:O13 created blank"
